$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": insert a new client row "JOWIN SA" for advisor
# "CASTRO ALCIVAR EDA MARIA" at row 30 (alphabetically before "MAD&DECO S.A.")
# shifting all subsequent rows down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(30).Insert()

$ws1.Range("A30").Value2 = "CASTRO ALCIVAR EDA MARIA"
$ws1.Range("B30").Value2 = "JOWIN SA"
$ws1.Range("C30:R30").Value2 = 0

# The final "N de 54" summary row (now shifted from row 56 to row 57) needs
# its denominator updated to reflect the new total row count (55).
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols1) {
    $cell = $ws1.Range($col + "57")
    $cell.Value2 = $cell.Value2.Replace("de 54", "de 55")
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same insertion of "JOWIN SA" row at row 30.
# The grand-total row (shifted from row 56 to row 57) keeps its same sums
# since the inserted row contributes only zeros.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(30).Insert()

$ws2.Range("A30").Value2 = "CASTRO ALCIVAR EDA MARIA"
$ws2.Range("B30").Value2 = "JOWIN SA"
$ws2.Range("C30:G30").Value2 = 0
